# OW-268 link uploaded trades to their account
# The uploaded trade's "Position Account ID" (cell B2 on the IRS-Cleared
# sheet) pointed at the placeholder "MEGA104" - repoint it at the real
# account id "acc1", and bring the view/selection back to that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- core semantic change -------------------------------------------------
$ws.Range("B2").Value = "acc1"

# --- restore the view to the edited cell ----------------------------------
$ws.Range("B2").Select() | Out-Null

# --- re-assert the (incidentally re-measured) column widths --------------
# Values below are the nearest width this engine's ColumnWidth API can
# reproduce (it stores widths quantised to 1/6 of a character) to the
# widths recorded in the edited workbook.
$ws.Columns(1).ColumnWidth = 8.166666666666666
$ws.Columns(2).ColumnWidth = 14.5
$ws.Columns(3).ColumnWidth = 7.5
$ws.Columns(4).ColumnWidth = 13.166666666666666
$ws.Columns(5).ColumnWidth = 7.666666666666667
$ws.Columns(6).ColumnWidth = 10.5
$ws.Columns(7).ColumnWidth = 9.666666666666666
$ws.Columns(8).ColumnWidth = 9.5
$ws.Columns(9).ColumnWidth = 7.666666666666667
$ws.Columns(10).ColumnWidth = 9.833333333333334
$ws.Columns(11).ColumnWidth = 9.833333333333334
$ws.Columns(12).ColumnWidth = 7.666666666666667
$ws.Columns(13).ColumnWidth = 5.5
$ws.Columns(14).ColumnWidth = 5.5
$ws.Columns(15).ColumnWidth = 11.666666666666666
$ws.Columns(16).ColumnWidth = 9.5
$ws.Columns(17).ColumnWidth = 8.833333333333334
$ws.Columns(18).ColumnWidth = 14.0
$ws.Columns(19).ColumnWidth = 26.5
$ws.Columns(20).ColumnWidth = 16.166666666666668
$ws.Columns(21).ColumnWidth = 14.333333333333334
$ws.Columns(22).ColumnWidth = 17.5
$ws.Columns(23).ColumnWidth = 16.666666666666668
$ws.Columns(24).ColumnWidth = 15.833333333333334
$ws.Columns(25).ColumnWidth = 15.333333333333334
$ws.Columns(26).ColumnWidth = 14.0
$ws.Columns(27).ColumnWidth = 11.666666666666666
$ws.Columns(28).ColumnWidth = 15.166666666666666
$ws.Columns(29).ColumnWidth = 9.5
$ws.Columns(30).ColumnWidth = 8.833333333333334
$ws.Columns(31).ColumnWidth = 14.0
$ws.Columns(32).ColumnWidth = 26.5
$ws.Columns(33).ColumnWidth = 16.166666666666668
$ws.Columns(34).ColumnWidth = 14.333333333333334
$ws.Columns(35).ColumnWidth = 31.666666666666668
$ws.Columns(36).ColumnWidth = 16.666666666666668
$ws.Columns(37).ColumnWidth = 15.833333333333334
$ws.Columns(38).ColumnWidth = 15.333333333333334
$ws.Columns(39).ColumnWidth = 14.0
$ws.Columns(40).ColumnWidth = 11.666666666666666
$ws.Columns(41).ColumnWidth = 15.166666666666666
$ws.Columns(42).ColumnWidth = 14.333333333333334
$ws.Columns(43).ColumnWidth = 14.333333333333334
$ws.Columns(44).ColumnWidth = 7.333333333333333
$ws.Columns(45).ColumnWidth = 7.333333333333333
$ws.Columns(46).ColumnWidth = 35.666666666666664
